$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add Sheet2 and Sheet3 right after Sheet1 (so ordering is Sheet1, Sheet2, Sheet3)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

$labels = @("Overall Score", "Productivity", "Creativity", "Responsiveness")

$sheets = @($ws1, $ws2, $ws3)
$values = @(
    @("1685`n", "1663`n", "1736`n", "1603`n"),
    @("1689`n", "1675`n", "1719`n", "1641`n"),
    @("1681`n", "1664`n", "1701`n", "1669`n")
)

for ($s = 0; $s -lt $sheets.Length; $s++) {
    $ws = $sheets[$s]

    # Column A: labels
    for ($r = 1; $r -le 4; $r++) {
        $ws.Range("A$r").Value = $labels[$r - 1]
    }

    # Column B: benchmark scores, stored as literal text (trailing newline, like the source data)
    $ws.Range("B1:B4").NumberFormat = "@"
    for ($r = 1; $r -le 4; $r++) {
        $ws.Range("B$r").Value = $values[$s][$r - 1]
    }
    # Drop back to the default (unstyled) cell format now that the text is committed
    $ws.Range("B1:B4").Style = "Normal"

    # Column A is wide enough to show the longest label in full
    $ws.Columns.Item(1).ColumnWidth = 19.86

    # Clear the auto row-height bump the embedded newline triggers
    for ($r = 1; $r -le 4; $r++) {
        $ws.Rows.Item($r).AutoFit()
    }
}

# Restore Sheet1 as the active/selected sheet
$ws1.Activate()
$excel.ActiveWindow.SelectedSheets.Item(1).Select()
